$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Current layout (before edit):
#   row 7  -> ALL-VENT SYRUP 125ML
#   row 8  -> BI ALCOFAN 150 MG 30 TABS.
#   row 9  -> LAMIFEN 1% CREAM 15 GM
#   row 10 -> PANADOL ADVANCE 500 MG 48 TABLETS
#   row 11 -> مجموعه برد
#   row 12 -> totals row (P12)
#   row 13 -> footer row (date / page / developer credit)
#
# Target layout (after edit):
#   row 7  -> ALL-VENT SYRUP 125ML            (unchanged)
#   row 8  -> BI ALCOFAN 150 MG 30 TABS.       (unchanged)
#   row 9  -> CETAL 500MG 20 TAB               (NEW)
#   row 10 -> IVYROSPAN SYRUP 100 ML           (NEW)
#   row 11 -> LAMIFEN 1% CREAM 15 GM           (shifted down, was row 9)
#   row 12 -> PANADOL ADVANCE 500 MG 48 TABLETS(shifted down, was row 10)
#   row 13 -> مجموعه برد                       (shifted down, was row 11)
#   row 14 -> معجون سيجنال 50 مل                (NEW)
#   row 15 -> totals row (P15, new sum)
#   row 16 -> footer row (updated timestamp)
# ---------------------------------------------------------------------------

# Insert two fresh rows at position 9 (for CETAL + IVYROSPAN); this pushes the
# old rows 9-13 down to 11-15.
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

# Insert one more fresh row right before the totals row (now at row 14) for
# the new معجون سيجنال line; pushes totals/footer down to 15/16.
$ws.Rows.Item(14).Insert()

# ---------------------------------------------------------------------------
# Re-apply cell formatting (number format / font / borders / alignment) to the
# three freshly-inserted, still-blank rows by pulling it from an existing
# item row that already carries the right look.
# ---------------------------------------------------------------------------
$ws.Range("A7:Q7").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)

$ws.Range("A8:Q8").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)

$ws.Range("A13:Q13").Copy()
$ws.Range("A14:Q14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row heights - match the target exactly.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 24.75
$ws.Rows.Item(14).RowHeight = 25.5
$ws.Rows.Item(15).RowHeight = 24.75
$ws.Rows.Item(16).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Merge the cell groups for the three new rows, matching the established
# A:B / C:G / H:K / L:M / N:O pattern used by every item row.
# ---------------------------------------------------------------------------
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

# ---------------------------------------------------------------------------
# Values - item rows 7 & 8 stay as-is. Fill in the rest top to bottom.
# ---------------------------------------------------------------------------

# Row 9 : CETAL 500MG 20 TAB
$ws.Cells.Item(9, 1).Value = 3
$ws.Cells.Item(9, 3).Value = "CETAL 500MG 20 TAB"
$ws.Cells.Item(9, 8).Value = "0:1"
$ws.Cells.Item(9, 12).Value = "1"
$ws.Cells.Item(9, 14).Value = "24.00"
$ws.Cells.Item(9, 16).Value = "12.0000"
$ws.Cells.Item(9, 17).Value = "0:1"

# Row 10 : IVYROSPAN SYRUP 100 ML
$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 3).Value = "IVYROSPAN SYRUP 100 ML"
$ws.Cells.Item(10, 8).Value = "0:0"
$ws.Cells.Item(10, 12).Value = "1"
$ws.Cells.Item(10, 14).Value = "55.00"
$ws.Cells.Item(10, 16).Value = "55.0000"
$ws.Cells.Item(10, 17).Value = "1:0"

# Row 11 : LAMIFEN 1% CREAM 15 GM (pre-existing data, item number stays 5)
$ws.Cells.Item(11, 1).Value = 5
$ws.Cells.Item(11, 3).Value = "LAMIFEN 1% CREAM 15 GM"
$ws.Cells.Item(11, 8).Value = "3:0"
$ws.Cells.Item(11, 12).Value = "1"
$ws.Cells.Item(11, 14).Value = "18.00"
$ws.Cells.Item(11, 16).Value = "18.0000"
$ws.Cells.Item(11, 17).Value = "1:0"

# Row 12 : PANADOL ADVANCE 500 MG 48 TABLETS
$ws.Cells.Item(12, 1).Value = 6
$ws.Cells.Item(12, 3).Value = "PANADOL ADVANCE 500 MG 48 TABLETS"
$ws.Cells.Item(12, 8).Value = "5:3"
$ws.Cells.Item(12, 12).Value = "1"
$ws.Cells.Item(12, 14).Value = "92.00"
$ws.Cells.Item(12, 16).Value = "23.0000"
$ws.Cells.Item(12, 17).Value = "0:1"

# Row 13 : مجموعه برد
$ws.Cells.Item(13, 1).Value = 7
$ws.Cells.Item(13, 3).Value = "مجموعه برد"
$ws.Cells.Item(13, 8).Value = "0:0"
$ws.Cells.Item(13, 12).Value = "0"
$ws.Cells.Item(13, 14).Value = "8.00"
$ws.Cells.Item(13, 16).Value = "8.0000"
$ws.Cells.Item(13, 17).Value = "1:0"

# Row 14 : معجون سيجنال 50 مل
$ws.Cells.Item(14, 1).Value = 8
$ws.Cells.Item(14, 3).Value = "معجون سيجنال 50 مل"
$ws.Cells.Item(14, 8).Value = "14:0"
$ws.Cells.Item(14, 12).Value = "0"
$ws.Cells.Item(14, 14).Value = "35.00"
$ws.Cells.Item(14, 16).Value = "35.0000"
$ws.Cells.Item(14, 17).Value = "1:0"

# Row 15 : totals
$ws.Cells.Item(15, 16).Value = 211.73

# Row 16 : footer (updated run timestamp)
$ws.Cells.Item(16, 1).Value = "Saturday, 31 May, 2025 10:00 AM"
$ws.Cells.Item(16, 7).Value = "1/1"
$ws.Cells.Item(16, 11).Value = "developed by : Abdelaziz Talaat"
